$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 2347.9
$ws.Range("I32").Value = 1400
$ws.Range("J32").Value = 2754.1428
$ws.Range("K32").Value = 1400
$ws.Range("L32").Value = 2754.1428
$ws.Range("M32").Value = -1074
$ws.Range("N32").Value = -3406.1428
$ws.Range("H45").Value = 193.5
$ws.Range("I45").Value = 187
$ws.Range("J45").Value = 200
$ws.Range("K45").Value = 561
$ws.Range("L45").Value = 600
$ws.Range("M45").Value = -369
$ws.Range("N45").Value = -984
$ws.Range("H125").Value = 350
$ws.Range("I125").Value = 300.33334
$ws.Range("J125").Value = 499
$ws.Range("K125").Value = 2703.00006
$ws.Range("L125").Value = 4491
$ws.Range("M125").Value = -243.0000600000003
$ws.Range("N125").Value = -9411
$ws.Range("H137").Value = 1236.75
$ws.Range("I137").Value = 1210.875
$ws.Range("J137").Value = 1288.5
$ws.Range("K137").Value = 3632.625
$ws.Range("L137").Value = 3865.5
$ws.Range("M137").Value = -1082.625
$ws.Range("N137").Value = -8965.5
$ws.Range("H138").Value = 1363.92
$ws.Range("I138").Value = 693.05884
$ws.Range("J138").Value = 1709.5151
$ws.Range("K138").Value = 2079.17652
$ws.Range("L138").Value = 5128.5453
$ws.Range("M138").Value = 3060.82348
$ws.Range("N138").Value = -15408.5453

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 11841.556
$ws.Range("I2").Value = 640.25
$ws.Range("J2").Value = 20802.6
$ws.Range("K2").Value = 640.25
$ws.Range("L2").Value = 20802.6
$ws.Range("M2").Value = -527.25
$ws.Range("N2").Value = -21028.6
$ws.Range("H4").Value = 236.75
$ws.Range("I4").Value = 282.33334
$ws.Range("K4").Value = 282.33334
$ws.Range("M4").Value = -166.33334
$ws.Range("H32").Value = 5103.931
$ws.Range("I32").Value = 5259.778
$ws.Range("J32").Value = 3000
$ws.Range("K32").Value = 5259.778
$ws.Range("L32").Value = 3000
$ws.Range("M32").Value = -4972.778
$ws.Range("N32").Value = -3574
$ws.Range("H61").Value = 1194.8334
$ws.Range("I61").Value = 939.8182
$ws.Range("J61").Value = 4000
$ws.Range("K61").Value = 939.8182
$ws.Range("L61").Value = 4000
$ws.Range("M61").Value = -727.8182
$ws.Range("N61").Value = -4424
$ws.Range("H62").Value = 75000
$ws.Range("J62").Value = 75000
$ws.Range("L62").Value = 75000
$ws.Range("N62").Value = -76248
$ws.Range("H65").Value = 75000
$ws.Range("J65").Value = 75000
$ws.Range("L65").Value = 225000
$ws.Range("N65").Value = -231240
$ws.Range("H110").Value = 2193.6
$ws.Range("I110").Value = 1753.3334
$ws.Range("K110").Value = 1753.3334
$ws.Range("M110").Value = 291.6666
$ws.Range("H116").Value = 11841.556
$ws.Range("I116").Value = 640.25
$ws.Range("J116").Value = 20802.6
$ws.Range("K116").Value = 640.25
$ws.Range("L116").Value = 20802.6
$ws.Range("M116").Value = 1653.75
$ws.Range("N116").Value = -25390.6
$ws.Range("H128").Value = 99990
$ws.Range("J128").Value = 99990
$ws.Range("L128").Value = 99990
$ws.Range("N128").Value = -109950
$ws.Range("H132").Value = 2197.6365
$ws.Range("I132").Value = 1860.5264
$ws.Range("K132").Value = 5581.5792
$ws.Range("M132").Value = -3051.5792
$ws.Range("H136").Value = 1194.8334
$ws.Range("I136").Value = 939.8182
$ws.Range("J136").Value = 4000
$ws.Range("K136").Value = 2819.4546
$ws.Range("L136").Value = 12000
$ws.Range("M136").Value = -269.4546
$ws.Range("N136").Value = -17100

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 11841.556
$ws.Range("I3").Value = 640.25
$ws.Range("J3").Value = 20802.6
$ws.Range("K3").Value = 640.25
$ws.Range("L3").Value = 20802.6
$ws.Range("M3").Value = -526.25
$ws.Range("N3").Value = -21030.6
$ws.Range("H86").Value = 3549.4614
$ws.Range("I86").Value = 3806.45
$ws.Range("J86").Value = 2692.8333
$ws.Range("K86").Value = 3806.45
$ws.Range("L86").Value = 2692.8333
$ws.Range("M86").Value = -2683.45
$ws.Range("N86").Value = -4938.8333
$ws.Range("H89").Value = 3549.4614
$ws.Range("I89").Value = 3806.45
$ws.Range("J89").Value = 2692.8333
$ws.Range("K89").Value = 19032.25
$ws.Range("L89").Value = 13464.1665
$ws.Range("M89").Value = -13416.25
$ws.Range("N89").Value = -24696.1665
$ws.Range("H105").Value = 83335010
$ws.Range("I105").Value = 125001820
$ws.Range("J105").Value = 1375.75
$ws.Range("K105").Value = 125001820
$ws.Range("L105").Value = 1375.75
$ws.Range("M105").Value = -125000073
$ws.Range("N105").Value = -4869.75
$ws.Range("H107").Value = 2077.4285
$ws.Range("I107").Value = 1558.875
$ws.Range("K107").Value = 1558.875
$ws.Range("M107").Value = 361.125
$ws.Range("H134").Value = 9230.77
$ws.Range("I134").Value = 1625.1
$ws.Range("J134").Value = 34583
$ws.Range("K134").Value = 4875.299999999999
$ws.Range("L134").Value = 103749
$ws.Range("M134").Value = -2340.299999999999
$ws.Range("N134").Value = -108819

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 681
$ws.Range("I107").Value = 492.55554
$ws.Range("J107").Value = 835.1818
$ws.Range("K107").Value = 492.55554
$ws.Range("L107").Value = 835.1818
$ws.Range("M107").Value = 1427.44446
$ws.Range("N107").Value = -4675.1818
$ws.Range("H132").Value = 2827.6155
$ws.Range("I132").Value = 2084.5557
$ws.Range("K132").Value = 6253.6671
$ws.Range("M132").Value = -3723.6671

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 100
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 100
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 600
$ws.Range("N2").Value = -826
$ws.Range("H106").Value = 3467.7273
$ws.Range("J106").Value = 3467.7273
$ws.Range("L106").Value = 10403.1819
$ws.Range("N106").Value = -12295.1819
$ws.Range("H107").Value = 531.125
$ws.Range("I107").Value = 111
$ws.Range("J107").Value = 628.0769
$ws.Range("K107").Value = 333
$ws.Range("L107").Value = 1884.2307
$ws.Range("M107").Value = 1587
$ws.Range("N107").Value = -5724.2307
$ws.Range("H132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("H137").Value = 14024
$ws.Range("I137").Value = 3092
$ws.Range("J137").Value = 27689
$ws.Range("K137").Value = 9276
$ws.Range("L137").Value = 83067
$ws.Range("M137").Value = -4176
$ws.Range("N137").Value = -93267
$ws.Range("M2").ClearContents()
$ws.Range("N132").ClearContents()

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4984.3335
$ws.Range("J80").Value = 4984.3335
$ws.Range("L80").Value = 4984.3335
$ws.Range("N80").Value = -6980.3335
$ws.Range("H83").Value = 4984.3335
$ws.Range("J83").Value = 4984.3335
$ws.Range("L83").Value = 24921.6675
$ws.Range("N83").Value = -34905.6675
$ws.Range("H113").Value = 2958.7368
$ws.Range("I113").Value = 1439.125
$ws.Range("J113").Value = 4063.9092
$ws.Range("K113").Value = 1439.125
$ws.Range("L113").Value = 4063.9092
$ws.Range("M113").Value = 730.875
$ws.Range("N113").Value = -8403.9092
$ws.Range("H122").Value = 752250
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("H132").Value = 2062.889
$ws.Range("I132").Value = 1608.9333
$ws.Range("K132").Value = 4826.7999
$ws.Range("M132").Value = -2296.7999
$ws.Range("M122").ClearContents()

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1685
$ws.Range("I16").Value = 1702.8889
$ws.Range("J16").Value = 1644.75
$ws.Range("K16").Value = 1702.8889
$ws.Range("L16").Value = 1644.75
$ws.Range("M16").Value = -1532.8889
$ws.Range("N16").Value = -1984.75
$ws.Range("H40").Value = 4333
$ws.Range("I40").Value = 4333
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 4333
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -4197
$ws.Range("H46").Value = 5948.5
$ws.Range("I46").Value = 3061.6667
$ws.Range("K46").Value = 3061.6667
$ws.Range("M46").Value = -2873.6667
$ws.Range("H61").Value = 1707.7273
$ws.Range("I61").Value = 1668.5714
$ws.Range("J61").Value = 1776.25
$ws.Range("K61").Value = 1668.5714
$ws.Range("L61").Value = 1776.25
$ws.Range("M61").Value = -1466.5714
$ws.Range("N61").Value = -2180.25
$ws.Range("H100").Value = 2175
$ws.Range("I100").Value = 1908.3334
$ws.Range("J100").Value = 2975
$ws.Range("K100").Value = 1908.3334
$ws.Range("L100").Value = 2975
$ws.Range("M100").Value = -1367.3334
$ws.Range("N100").Value = -4057
$ws.Range("H113").Value = 1707.7273
$ws.Range("I113").Value = 1668.5714
$ws.Range("J113").Value = 1776.25
$ws.Range("K113").Value = 1668.5714
$ws.Range("L113").Value = 1776.25
$ws.Range("M113").Value = 501.4286
$ws.Range("N113").Value = -6116.25
$ws.Range("H132").Value = 85371.75
$ws.Range("I132").Value = 1182.25
$ws.Range("J132").Value = 253750.75
$ws.Range("K132").Value = 3546.75
$ws.Range("L132").Value = 761252.25
$ws.Range("M132").Value = -1016.75
$ws.Range("N132").Value = -766312.25
$ws.Range("H136").Value = 2182.3333
$ws.Range("I136").Value = 1936.8
$ws.Range("J136").Value = 2357.7144
$ws.Range("K136").Value = 5810.4
$ws.Range("L136").Value = 7073.1432
$ws.Range("M136").Value = -3260.4
$ws.Range("N136").Value = -12173.1432
$ws.Range("N40").ClearContents()

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H58").Value = 15042
$ws.Range("I58").Value = 15042
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 15042
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = -14734
$ws.Range("H113").Value = 918.25
$ws.Range("I113").Value = 389.2
$ws.Range("J113").Value = 1800
$ws.Range("K113").Value = 1167.6
$ws.Range("L113").Value = 5400
$ws.Range("M113").Value = 1002.4
$ws.Range("N113").Value = -9740
$ws.Range("H126").Value = 125000880
$ws.Range("I126").Value = 166667650
$ws.Range("J126").Value = 609
$ws.Range("K126").Value = 500002950
$ws.Range("L126").Value = 1827
$ws.Range("M126").Value = -500000480
$ws.Range("N126").Value = -6767
$ws.Range("H132").Value = 3951
$ws.Range("I132").Value = 3576.8572
$ws.Range("J132").Value = 4998.6
$ws.Range("K132").Value = 10730.5716
$ws.Range("L132").Value = 14995.8
$ws.Range("M132").Value = -8200.571599999999
$ws.Range("N132").Value = -20055.8
$ws.Range("H136").Value = 1151.5
$ws.Range("I136").Value = 1151.5
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 3454.5
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -904.5
$ws.Range("N58").ClearContents()
$ws.Range("N136").ClearContents()
